# "New crime data collected" - weekly CompStat refresh (10th Precinct,
# report week 9/8/2025 - 9/14/2025, Volume 32 Number 37).
#
# Updates the two header strings (issue number + reporting week dates) and
# refreshes every statistic cell in the crime-complaints table (rows 15-31,
# columns C:N) to the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header text: masthead issue number and the "week covering" dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 32   Number  37"
$ws.Range("C9").Value  = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# ---------------------------------------------------------------------
# 2) Bulk numeric refresh for the Crime Complaints grid (rows 15-31).
#    Cells that keep their existing type (number-in/number-out) are just
#    overwritten in place; the handful of cells that flip between a
#    numeric value and the "no data" placeholder text are handled
#    separately below so their style/type lands correctly.
# ---------------------------------------------------------------------
$updates = @{
    "N15" = -36.842105263157

    "C16" = 2
    "E16" = 0
    "F16" = 8
    "G16" = 10
    "H16" = -20
    "I16" = 77
    "J16" = 76
    "K16" = 1.315789473684
    "L16" = -22.222222222222
    "M16" = 10
    "N16" = -82.579185520362

    "C17" = 2
    "D17" = 1
    "E17" = 100
    "F17" = 10
    "G17" = 10
    "I17" = 108
    "J17" = 88
    "K17" = 22.727272727272
    "L17" = -7.692307692307
    "M17" = 14.893617021276
    "N17" = -39.664804469273

    "D18" = 1
    "E18" = 100
    "F18" = 10
    "G18" = 7
    "H18" = 42.857142857142
    "I18" = 78
    "J18" = 72
    "K18" = 8.333333333333
    "L18" = -18.75
    "M18" = 6.849315068493
    "N18" = -81.560283687943

    "C19" = 17
    "D19" = 17
    "F19" = 53
    "G19" = 62
    "H19" = -14.516129032258
    "I19" = 431
    "J19" = 495
    "K19" = -12.929292929292
    "L19" = -21.636363636363
    "M19" = 0.700934579439
    "N19" = -24.253075571177

    "D20" = 4
    "E20" = -100
    "F20" = 3
    "G20" = 10
    "H20" = -70
    "J20" = 38
    "K20" = -42.105263157894
    "L20" = -62.068965517241
    "N20" = -93.604651162790

    "D21" = 25
    "E21" = -8
    "F21" = 84
    "H21" = -15.151515151515
    "I21" = 729
    "J21" = 779
    "K21" = -6.418485237483
    "L21" = -21.103896103896
    "M21" = 3.551136363636
    "N21" = -63.200403836446

    "I22" = 22
    "K22" = 0
    "L22" = -12
    "M22" = 144.444444444444

    "F23" = 3
    "G23" = 1
    "H23" = 200
    "I23" = 29
    "K23" = -3.333333333333
    "L23" = -17.142857142857
    "M23" = -12.121212121212

    "C24" = 22
    "D24" = 13
    "E24" = 69.230769230769
    "F24" = 96
    "G24" = 60
    "H24" = 60
    "I24" = 629
    "J24" = 480
    "K24" = 31.041666666666
    "L24" = 7.521367521367
    "M24" = 1.779935275080

    "C25" = 14
    "D25" = 7
    "E25" = 100
    "F25" = 53
    "H25" = 60.606060606060
    "I25" = 337
    "J25" = 212
    "K25" = 58.962264150943
    "L25" = -2.034883720930

    "C26" = 10
    "D26" = 6
    "E26" = 66.666666666666
    "F26" = 24
    "G26" = 25
    "H26" = -4
    "I26" = 215
    "J26" = 224
    "K26" = -4.017857142857
    "L26" = -10.788381742738
    "M26" = -10.041841004184

    "D28" = 2
    "E28" = -50
    "F28" = 3
    "G28" = 6
    "I28" = 38
    "J28" = 43
    "K28" = -11.627906976744
    "L28" = -13.636363636363

    "J31" = 6
    "K31" = -33.333333333333
    "L31" = -73.333333333333
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---------------------------------------------------------------------
# 3) Cells that flip type this week between a real number and the
#    "no data collected" placeholder text ("0" / "***.*"). Assigning the
#    placeholder text while the cell is in its normal numeric format
#    would just coerce "0" back to a number, so the format is forced to
#    Text first; likewise, turning a placeholder back into a number
#    needs its numeric format restored. Re-pasting the format (only)
#    from a donor cell that already carries the right style keeps the
#    style id identical to the ones used elsewhere in the sheet for the
#    same role (General/text vs "#,##0" vs "#,##0.0").
# ---------------------------------------------------------------------

function Set-PlaceholderText($ref, $text, $formatDonor) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($formatDonor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

function Set-NumberFromPlaceholder($ref, $value, $formatDonor) {
    $ws.Range($ref).Value = $value
    $ws.Range($formatDonor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

# Row 20 - G.L.A.: Week-to-date 2025 complaints not yet collected.
Set-PlaceholderText "C20" "0" "C14"

# Row 22 - Transit: complaints swap which year has the "no data" figure.
Set-NumberFromPlaceholder "C22" 1 "D16"
Set-PlaceholderText       "D22" "0" "C14"
Set-PlaceholderText       "E22" "***.*" "C14"

# Row 31 - Hate Crimes: 28-day window now has real counts on both sides.
Set-NumberFromPlaceholder "D31" 1    "D16"
Set-NumberFromPlaceholder "E31" -100 "E16"
Set-NumberFromPlaceholder "G31" 1    "D16"
Set-NumberFromPlaceholder "H31" -100 "E16"
